# GenEditID submission spreadsheet: rename the "DesireEditedSequences" tab
# to "TargetedSearch", and rename "Documentation" to "Help".
# (python/tests/pytest_min.xlsx — see commit message: "rename tab in
# submission spreadsheet from DesireEditedSequences to TargetedSearch")

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Documentation").Name = "Help"
$wb.Worksheets.Item("DesireEditedSequences").Name = "TargetedSearch"

# The Layout sheet's saved cursor/selection moves from the whole of
# column C down to just C2.
$null = $wb.Worksheets.Item("Layout").Range("C2").Select()

# The active/selected workbook tab moves off the Amplicon sheet and
# onto the newly-renamed Help sheet.
$null = $wb.Worksheets.Item("Help").Select()

Write-Output "renamed Documentation -> Help, DesireEditedSequences -> TargetedSearch"
